{"js": "// \"fix linebreaks in headers\"\n//\n// 1. Remove the stray \"_GoBack\" bookmark that sits after the \"center\"\n//    paragraph (left over from a previous edit session).\n// 2. Drop the trailing space in the \"Heading 3 \" paragraph so it reads\n//    \"Heading 3\".\n// 3. Re-insert the \"_GoBack\" bookmark where Word would actually have\n//    left it after the last real edit: in the middle of the \"Heading 6\"\n//    text (between \"Headin\" and \"g 6\") - this is what naturally causes\n//    the bookmark id renumbering seen in the diff.\n// 4. Stop forcing \"keep with next paragraph\" on Heading 1-6 so a heading\n//    followed by a long paragraph no longer forces an awkward\n//    line/page break.\n\n// --- 1. remove the old _GoBack bookmark -----------------------------\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\n// --- 2. trim the trailing space from \"Heading 3 \" --------------------\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nlet heading3 = null;\nlet heading6 = null;\nfor (const p of paragraphs.items) {\n  if (p.text === \"Heading 3 \") {\n    heading3 = p;\n  } else if (p.text === \"Heading 6\") {\n    heading6 = p;\n  }\n}\n\nif (heading3) {\n  const spaces = heading3.search(\" \", { matchCase: true });\n  spaces.load(\"items\");\n  await context.sync();\n  // the trailing space is the last match inside the paragraph\n  const trailingSpace = spaces.items[spaces.items.length - 1];\n  trailingSpace.delete();\n  await context.sync();\n}\n\n// --- 3. re-insert _GoBack inside \"Heading 6\", splitting the run ------\nif (heading6) {\n  const tail = heading6.search(\"g 6\");\n  tail.load(\"items\");\n  await context.sync();\n  const splitPoint = tail.items[0].getRange(\"Start\");\n  splitPoint.insertBookmark(\"_GoBack\");\n  await context.sync();\n}\n\n// --- 4. stop keeping Heading 1-6 glued to the following paragraph ----\nconst headingStyleNames = [\n  \"Heading 1\",\n  \"Heading 2\",\n  \"Heading 3\",\n  \"Heading 4\",\n  \"Heading 5\",\n  \"Heading 6\",\n];\nfor (const name of headingStyleNames) {\n  const style = context.document.getStyles().getByNameOrNullObject(name);\n  style.paragraphFormat.keepWithNext = false;\n}\nawait context.sync();\n", "ps1": "# \"fix linebreaks in headers\"\n#\n# 1. Remove the stray \"_GoBack\" bookmark that sits after the \"center\"\n#    paragraph (left over from a previous edit session).\n# 2. Drop the trailing space in the \"Heading 3 \" paragraph so it reads\n#    \"Heading 3\".\n# 3. Re-insert the \"_GoBack\" bookmark where Word would actually have\n#    left it after the last real edit: in the middle of the \"Heading 6\"\n#    text (between \"Headin\" and \"g 6\") - this is what naturally causes\n#    the bookmark id renumbering seen in the diff.\n# 4. Stop forcing \"keep with next paragraph\" on Heading 1-6 so a heading\n#    followed by a long paragraph no longer forces an awkward\n#    line/page break.\n\n$d = $word.ActiveDocument\n\n# --- 1. remove the old _GoBack bookmark ------------------------------\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n\n# --- 2. trim the trailing space from \"Heading 3 \" ---------------------\n$find = $d.Content\n$find.Find.Execute(\"Heading 3 \")\nif ($find.Find.Found) {\n    $trailing = $d.Range($find.End - 1, $find.End)\n    $trailing.Delete()\n}\n\n# --- 3. re-insert _GoBack inside \"Heading 6\", splitting the run -------\n$h6 = $d.Content\n$h6.Find.Execute(\"g 6\")\nif ($h6.Find.Found) {\n    $splitPoint = $d.Range($h6.Start, $h6.Start)\n    $d.Bookmarks.Add(\"_GoBack\", $splitPoint) | Out-Null\n}\n\n# --- 4. stop keeping Heading 1-6 glued to the following paragraph -----\n$headingStyleNames = @(\"Heading 1\", \"Heading 2\", \"Heading 3\", \"Heading 4\", \"Heading 5\", \"Heading 6\")\nforeach ($name in $headingStyleNames) {\n    $style = $d.Styles($name)\n    $style.ParagraphFormat.KeepWithNext = $False\n}\n"}
